$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.138.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "'2.468.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'491.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'153.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.16%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").Value = "'2.474.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +6.28%  "
$ws.Range("D11").Value = "'5.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "'2.896.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "'57.209.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.14%  "
$ws.Range("D16").Value = "'21.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "'2.478.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "'4.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.29%  "
$ws.Range("D20").Value = "'326.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.01%  "
$ws.Range("D21").Value = "'10.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'5.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'0.164"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "'2.559.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'7.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("D30").Value = "'0.0₃0821"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.14%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'150.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  +5.16%  "
$ws.Range("D34").Value = "'18.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'5.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").Value = "'0.903"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.17%  "
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").Value = "'3.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.98%  "
$ws.Range("D39").Value = "'1.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.91%  "
$ws.Range("D40").Value = "'34.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "'0.0560"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'0.610"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0961"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.65%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("D47").Value = "'267.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.30%  "
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").Value = "'17.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.97%  "
$ws.Range("E51").Value = "  +28.20%  "
